# Add two new columns (I: "I0", J: "IF") to the sheet, mirroring the
# existing header style (bold/centered/bordered, same as H1) and filling
# in the per-row numeric data for rows 2-21.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Headers -----------------------------------------------------------
# Copy H1's style onto I1/J1 so the new headers look like the rest of the
# header row, then overwrite the copied text with the new header labels.
$ws.Range("H1").Copy($ws.Range("I1"))
$ws.Range("H1").Copy($ws.Range("J1"))
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data ----------------------------------------------------------------
# row => (I value, J value)
$data = @{
    2  = @(2, 7)
    3  = @(1, 4)
    4  = @(1, 5)
    5  = @(1, 4)
    6  = @(1, 2)
    7  = @(1, 6)
    8  = @(1, 6)
    9  = @(1, 5)
    10 = @(1, 5)
    11 = @(1, 6)
    12 = @(4, 6)
    13 = @(3, 6)
    14 = @(6, 8)
    15 = @(1, 6)
    16 = @(1, 5)
    17 = @(5, 6)
    18 = @(9, 9)
    19 = @(1, 3)
    20 = @(4, 5)
    21 = @(3, 4)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]   # column I
    $ws.Cells.Item($row, 10).Value = $vals[1]  # column J
}
